$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers J1, K1 to the new "ExpectedRecordDuration_*" naming scheme
$ws.Range("J1").Value = "ExpectedRecordDuration_1"
$ws.Range("K1").Value = "ExpectedRecordDuration_2"

# Add new trailing columns S, T (ExpectedRecordDuration_3 / _4) with their row2 values
$ws.Range("S1").Value = "ExpectedRecordDuration_3"
$ws.Range("T1").Value = "ExpectedRecordDuration_4"
$ws.Range("S2").Value = 600
$ws.Range("T2").Value = 600

# Copy the number-as-text ("quote prefix") style from J2 onto the new S2 cell
$ws.Range("J2").Copy()
$ws.Range("S2").PasteSpecial(-4122)  # xlPasteFormats

# Rename L1 header from OmicronFile -> OmicronFile_1
$ws.Range("L1").Value = "OmicronFile_1"

# Add new trailing column U (OmicronFile_2) with its row2 value
$ws.Range("U1").Value = "OmicronFile_2"
$ws.Range("U2").Value = "CAM_734_1.seq"

# Copy the style from L2 onto the new U2 cell
$ws.Range("L2").Copy()
$ws.Range("U2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the selection to match the authored state
$ws.Range("L1").Select() | Out-Null
